$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 499.5
$ws.Range("I8").Value = 40.4
$ws.Range("J8").Value = 2795
$ws.Range("K8").Value = 121.2
$ws.Range("L8").Value = 8385
$ws.Range("M8").Value = 17.80000000000001
$ws.Range("N8").Value = -8663
$ws.Range("H32").Value = 4000
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H96").Value = 13572.375
$ws.Range("I96").Value = 17179.834
$ws.Range("J96").Value = 2750
$ws.Range("K96").Value = 51539.50199999999
$ws.Range("L96").Value = 8250
$ws.Range("M96").Value = -50166.50199999999
$ws.Range("N96").Value = -10996
$ws.Range("H98").Value = 1275.5
$ws.Range("I98").Value = 1118.6154
$ws.Range("K98").Value = 1118.6154
$ws.Range("M98").Value = 379.3846000000001
$ws.Range("H116").Value = 6616.636
$ws.Range("I116").Value = 4199.5
$ws.Range("J116").Value = 7997.857
$ws.Range("K116").Value = 4199.5
$ws.Range("L116").Value = 7997.857
$ws.Range("M116").Value = -757.5
$ws.Range("N116").Value = -14881.857
$ws.Range("H122").Value = 1275.5
$ws.Range("I122").Value = 1118.6154
$ws.Range("K122").Value = 3355.8462
$ws.Range("M122").Value = -905.8462
$ws.Range("H137").Value = 2484.8572
$ws.Range("I137").Value = 2297.5
$ws.Range("J137").Value = 2559.8
$ws.Range("K137").Value = 6892.5
$ws.Range("L137").Value = 7679.400000000001
$ws.Range("M137").Value = -4342.5
$ws.Range("N137").Value = -12779.4
$ws.Range("H138").Value = 5203.9375
$ws.Range("I138").Value = 4558
$ws.Range("J138").Value = 5849.875
$ws.Range("K138").Value = 13674
$ws.Range("L138").Value = 17549.625
$ws.Range("M138").Value = -8534
$ws.Range("N138").Value = -27829.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1996.3529
$ws.Range("I2").Value = 1272.3
$ws.Range("J2").Value = 3030.7144
$ws.Range("K2").Value = 1272.3
$ws.Range("L2").Value = 3030.7144
$ws.Range("M2").Value = -1159.3
$ws.Range("N2").Value = -3256.7144
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H116").Value = 1996.3529
$ws.Range("I116").Value = 1272.3
$ws.Range("J116").Value = 3030.7144
$ws.Range("K116").Value = 1272.3
$ws.Range("L116").Value = 3030.7144
$ws.Range("M116").Value = 1021.7
$ws.Range("N116").Value = -7618.7144

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1996.3529
$ws.Range("I3").Value = 1272.3
$ws.Range("J3").Value = 3030.7144
$ws.Range("K3").Value = 1272.3
$ws.Range("L3").Value = 3030.7144
$ws.Range("M3").Value = -1158.3
$ws.Range("N3").Value = -3258.7144
$ws.Range("H99").Value = 3625.6667
$ws.Range("I99").Value = 3688.75
$ws.Range("J99").Value = 3499.5
$ws.Range("K99").Value = 3688.75
$ws.Range("L99").Value = 3499.5
$ws.Range("M99").Value = -2190.75
$ws.Range("N99").Value = -6495.5
$ws.Range("H107").Value = 1303.6
$ws.Range("I107").Value = 526
$ws.Range("J107").Value = 4414
$ws.Range("K107").Value = 526
$ws.Range("L107").Value = 4414
$ws.Range("M107").Value = 1394
$ws.Range("N107").Value = -8254

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 189463.5
$ws.Range("I94").Value = 278889
$ws.Range("K94").Value = 278889
$ws.Range("M94").Value = -278438
$ws.Range("H134").Value = 2405.7856
$ws.Range("I134").Value = 2168.1
$ws.Range("K134").Value = 6504.299999999999
$ws.Range("M134").Value = -3969.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 1791.6364
$ws.Range("I38").Value = 2234.125
$ws.Range("J38").Value = 611.6667
$ws.Range("K38").Value = 6702.375
$ws.Range("L38").Value = 1835.0001
$ws.Range("M38").Value = -6355.375
$ws.Range("N38").Value = -2529.0001
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("N59").ClearContents()
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H107").Value = 4998.8
$ws.Range("I107").Value = 4998.3335
$ws.Range("J107").Value = 4999.5
$ws.Range("K107").Value = 14995.0005
$ws.Range("L107").Value = 14998.5
$ws.Range("M107").Value = -13075.0005
$ws.Range("N107").Value = -18838.5
$ws.Range("H137").Value = 2662.889
$ws.Range("I137").Value = 2261.6667
$ws.Range("J137").Value = 2863.5
$ws.Range("K137").Value = 6785.000100000001
$ws.Range("L137").Value = 8590.5
$ws.Range("M137").Value = -1685.000100000001
$ws.Range("N137").Value = -18790.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 17113.857
$ws.Range("I80").Value = 3166.6667
$ws.Range("J80").Value = 27574.25
$ws.Range("K80").Value = 3166.6667
$ws.Range("L80").Value = 27574.25
$ws.Range("M80").Value = -2168.6667
$ws.Range("N80").Value = -29570.25
$ws.Range("H83").Value = 17113.857
$ws.Range("I83").Value = 3166.6667
$ws.Range("J83").Value = 27574.25
$ws.Range("K83").Value = 15833.3335
$ws.Range("L83").Value = 137871.25
$ws.Range("M83").Value = -10841.3335
$ws.Range("N83").Value = -147855.25
$ws.Range("H102").Value = 1937.2858
$ws.Range("I102").Value = 2072.2
$ws.Range("J102").Value = 1600
$ws.Range("K102").Value = 2072.2
$ws.Range("L102").Value = 1600
$ws.Range("M102").Value = -450.1999999999998
$ws.Range("N102").Value = -4844

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3794.25
$ws.Range("I132").Value = 3792.6667
$ws.Range("K132").Value = 11378.0001
$ws.Range("M132").Value = -8848.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 16010.667
$ws.Range("J45").Value = 16010.667
$ws.Range("L45").Value = 16010.667
$ws.Range("N45").Value = -16992.667
$ws.Range("H81").Value = 1001000.3
$ws.Range("I81").Value = 1175.125
$ws.Range("J81").Value = 5000301
$ws.Range("K81").Value = 2350.25
$ws.Range("L81").Value = 10000602
$ws.Range("M81").Value = -1289.25
$ws.Range("N81").Value = -10002724
$ws.Range("H84").Value = 1001000.3
$ws.Range("I84").Value = 1175.125
$ws.Range("J84").Value = 5000301
$ws.Range("K84").Value = 11751.25
$ws.Range("L84").Value = 50003010
$ws.Range("M84").Value = -6447.25
$ws.Range("N84").Value = -50013618
$ws.Range("H126").Value = 1778.625
$ws.Range("I126").Value = 1747
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5241
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -2771
$ws.Range("N126").Value = -10940
